$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("M1", "Matn1", "Itgb1", "ECs", 3, 1, 0.2053186666666666, 0.615956, 0.420847658662702, 0.4208476586627021, 3, 1, 163.7119853333333, 491.135956, 0.2754003062401033, 0.2754003062401033, 33.61312654599288, 302.518138913936, 0.1159015740761386, 0.1159015740761386),
  @("M1", "Matn1", "Itgb1", "FAPs", 3, 1, 0.2053186666666666, 0.615956, 0.420847658662702, 0.4208476586627021, 3, 1, 172.558497, 517.675491, 0.290282124557779, 0.290282124557779, 35.429480526044, 318.8653247343959, 0.1221645524717761, 0.1221645524717761),
  @("M1", "Matn1", "Itgb1", "M1", 3, 1, 0.2053186666666666, 0.615956, 0.420847658662702, 0.4208476586627021, 3, 1, 111.4881643333333, 334.464493, 0.1875481171218523, 0.1875481171218523, 22.89060125003422, 206.015411250308, 0.07892918597732976, 0.07892918597732977),
  @("M1", "Matn1", "Itgb1", "M2", 3, 1, 0.2053186666666666, 0.615956, 0.420847658662702, 0.4208476586627021, 3, 1, 86.95798233333333, 260.873947, 0.1462828449356383, 0.1462828449356383, 17.85409698870355, 160.686872898332, 0.06156279279368249, 0.0615627927936825),
  @("M1", "Matn1", "Itgb1", "sCs", 3, 1, 0.2053186666666666, 0.615956, 0.420847658662702, 0.4208476586627021, 3, 1, 59.73436333333333, 179.20309, 0.100486607144627, 0.100486607144627, 12.26457983378222, 110.38121850404, 0.042289553343775, 0.04228955334377501),
  @("M2", "Matn1", "Itgb1", "ECs", 2, 0.6666666666666666, 0.2825506666666667, 0.847652, 0.579152341337298, 0.579152341337298, 3, 1, 163.7119853333333, 491.135956, 0.2754003062401033, 0.2754003062401033, 46.25693059725689, 416.3123753753119, 0.1594987321639647, 0.1594987321639647),
  @("M2", "Matn1", "Itgb1", "FAPs", 2, 0.6666666666666666, 0.2825506666666667, 0.847652, 0.579152341337298, 0.579152341337298, 3, 1, 172.558497, 517.675491, 0.290282124557779, 0.290282124557779, 48.756518366348, 438.8086652971319, 0.1681175720860029, 0.1681175720860029),
  @("M2", "Matn1", "Itgb1", "M1", 2, 0.6666666666666666, 0.2825506666666667, 0.847652, 0.579152341337298, 0.579152341337298, 3, 1, 111.4881643333333, 334.464493, 0.1875481171218523, 0.1875481171218523, 31.50105515782622, 283.509496420436, 0.1086189311445225, 0.1086189311445225),
  @("M2", "Matn1", "Itgb1", "M2", 2, 0.6666666666666666, 0.2825506666666667, 0.847652, 0.579152341337298, 0.579152341337298, 3, 1, 86.95798233333333, 260.873947, 0.1462828449356383, 0.1462828449356383, 24.57003588027156, 221.130322922444, 0.08472005214195585, 0.08472005214195585),
  @("M2", "Matn1", "Itgb1", "sCs", 2, 0.6666666666666666, 0.2825506666666667, 0.847652, 0.579152341337298, 0.579152341337298, 3, 1, 59.73436333333333, 179.20309, 0.100486607144627, 0.100486607144627, 16.87798418274222, 151.90185764468, 0.05819705380085197, 0.05819705380085197)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $rowData = $data[$i]
  for ($c = 0; $c -lt $rowData.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
  }
}
